# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Refresh case-count figures for several countries (new totals/new-cases/
#   active/recovered/critical/deaths).
# - Two pairs of rows swap order because their totals crossed after the
#   refresh (Bolivia/Oman and Islas Turcas y Caicos/Seychelles), so those
#   rows get the other country's name along with their own refreshed
#   numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 04:45"

# Row 30 - Kazajistan
$ws.Range("B30").Value = 92662
$ws.Range("C30").Value = 1069
$ws.Range("E30").Value = 30030

# Row 32 - China
$ws.Range("B32").Value = 84428
$ws.Range("C32").Value = 43
$ws.Range("D32").Value = 79013
$ws.Range("E32").Value = 781

# Row 34 - was Oman, now Bolivia (rows 34/35 swap order)
$ws.Range("A34").Value = "Bolivia"
$ws.Range("B34").Value = 80153
$ws.Range("C34").Value = 1360
$ws.Range("D34").Value = 24156
$ws.Range("E34").Value = 52844
$ws.Range("G34").Value = 89
$ws.Range("H34").Value = 3153

# Row 35 - was Bolivia, now Oman
$ws.Range("A35").Value = "Oman"
$ws.Range("B35").Value = 79159
$ws.Range("D35").Value = 61421
$ws.Range("E35").Value = 17317
$ws.Range("H35").Value = 421

# Row 39 - Belgica
$ws.Range("B39").Value = 69849
$ws.Range("C39").Value = 447
$ws.Range("D39").Value = 17590
$ws.Range("E39").Value = 42414

# Row 51 - Honduras
$ws.Range("B51").Value = 43197
$ws.Range("C51").Value = 512
$ws.Range("D51").Value = 5794
$ws.Range("E51").Value = 36026
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 1377

# Row 77 - Corea del Sur
$ws.Range("B77").Value = 14389
$ws.Range("C77").Value = 23
$ws.Range("D77").Value = 13280
$ws.Range("E77").Value = 808

# Row 152 - Jamaica
$ws.Range("B152").Value = 894
$ws.Range("C152").Value = 11
$ws.Range("E152").Value = 139
$ws.Range("G152").Value = 2
$ws.Range("H152").Value = 12

# Row 187 - was Seychelles, now Islas Turcas y Caicos (rows 187/188 swap order)
$ws.Range("A187").Value = "Islas Turcas y Caicos"
$ws.Range("B187").Value = 116
$ws.Range("C187").Value = 2
$ws.Range("D187").Value = 38
$ws.Range("E187").Value = 76
$ws.Range("H187").Value = 2

# Row 188 - was Islas Turcas y Caicos, now Seychelles
$ws.Range("A188").Value = "Seychelles"
$ws.Range("D188").Value = 39
$ws.Range("E188").Value = 75
$ws.Range("H188").Value = 0
